# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets to match the latest generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3327
$ws1.Range("F6").Value = 2128
$ws1.Range("F10").Value = 21
$ws1.Range("F11").Value = 1229
$ws1.Range("F13").Value = 1319

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3327
$ws4.Range("F6").Value = 2128
$ws4.Range("F11").Value = 21
$ws4.Range("F14").Value = 1229
$ws4.Range("F16").Value = 1319
